# Final version for the thesis.
# Adds a new "Other parameter values:" header in column H, widens that
# column, and tidies up the border formatting around the bottom of the
# "Mutation percentage" table (column A).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column header (H1) -------------------------------------------------
# Put the text in first, then copy the formatting that is already used by
# the other header cells (e.g. B1) so the new cell reuses the existing
# bold/centered/boxed header style instead of inventing a new one.
$ws.Range("H1").Value = "Other parameter values:"
$ws.Range("B1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# Widen the new column so the header text fits nicely.
$ws.Columns("H").ColumnWidth = 24.8

# --- Border tidy-up on column A ---------------------------------------------
# Drop the bottom border of A8 so the box around A8/A9 does not show a
# doubled line between the two cells.
$ws.Range("A8").Borders.Item(9).LineStyle = 0

# A10 is a blank spacer cell at the bottom of the table; remove its box
# border entirely while keeping its existing font/alignment formatting.
$ws.Range("A10").Borders.LineStyle = 0

# --- Reset the view back to the top-left cell -------------------------------
$ws.Range("A1").Select() | Out-Null
